# Auto-generated Excel COM-interop script
# Applies the numeric updates described by the commit diff to the
# Leve-profit sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 40
$ws.Range("H40").Value = 4000
$ws.Range("I40").Value = 0
$ws.Range("J40").Value = 4000
$ws.Range("K40").Value = 0
$ws.Range("L40").Value = 4000
$ws.Range("M40").ClearContents()
$ws.Range("N40").Value = -4350
# Row 42
$ws.Range("H42").Value = 279.23077
$ws.Range("I42").Value = 173.33333
$ws.Range("J42").Value = 370
$ws.Range("K42").Value = 519.99999
$ws.Range("L42").Value = 1110
$ws.Range("M42").Value = -289.99999
$ws.Range("N42").Value = -1570
# Row 43
$ws.Range("H43").Value = 812.86664
$ws.Range("I43").Value = 640.2
$ws.Range("J43").Value = 899.2
$ws.Range("K43").Value = 640.2
$ws.Range("L43").Value = 899.2
$ws.Range("M43").Value = -571.2
$ws.Range("N43").Value = -1037.2
# Row 86
$ws.Range("H86").Value = 1803.125
$ws.Range("J86").Value = 1810
$ws.Range("L86").Value = 1810
$ws.Range("N86").Value = -4056
# Row 89
$ws.Range("H89").Value = 1803.125
$ws.Range("J89").Value = 1810
$ws.Range("L89").Value = 9050
$ws.Range("N89").Value = -20282
# Row 106
$ws.Range("H106").Value = 2649.5833
$ws.Range("I106").Value = 2279
$ws.Range("J106").Value = 2914.2856
$ws.Range("K106").Value = 2279
$ws.Range("L106").Value = 2914.2856
$ws.Range("M106").Value = -1648
$ws.Range("N106").Value = -4176.2856
# Row 125
$ws.Range("H125").Value = 966.8182
$ws.Range("I125").Value = 737.1429
$ws.Range("J125").Value = 1074
$ws.Range("K125").Value = 6634.2861
$ws.Range("L125").Value = 9666
$ws.Range("M125").Value = -4174.2861
$ws.Range("N125").Value = -14586
# Row 132
$ws.Range("H132").Value = 3063.8076
$ws.Range("I132").Value = 2669.932
$ws.Range("J132").Value = 5230.125
$ws.Range("K132").Value = 8009.795999999999
$ws.Range("L132").Value = 15690.375
$ws.Range("M132").Value = -5479.795999999999
$ws.Range("N132").Value = -20750.375
# Row 141
$ws.Range("H141").Value = 1153.5
$ws.Range("I141").Value = 1126.8462
$ws.Range("J141").Value = 1500
$ws.Range("K141").Value = 3380.5386
$ws.Range("L141").Value = 4500
$ws.Range("M141").Value = 1799.4614
$ws.Range("N141").Value = -14860

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 13712811
$ws.Range("I32").Value = 17862378
$ws.Range("K32").Value = 17862378
$ws.Range("M32").Value = -17862091
# Row 101
$ws.Range("H101").Value = 80000
$ws.Range("J101").Value = 80000
$ws.Range("L101").Value = 80000
$ws.Range("N101").Value = -86490
# Row 102
$ws.Range("H102").Value = 3152.5
$ws.Range("I102").Value = 3499.5
$ws.Range("J102").Value = 2805.5
$ws.Range("K102").Value = 3499.5
$ws.Range("L102").Value = 2805.5
$ws.Range("M102").Value = -1877.5
$ws.Range("N102").Value = -6049.5
# Row 110
$ws.Range("H110").Value = 2494.2942
$ws.Range("I110").Value = 1649
$ws.Range("J110").Value = 3701.8572
$ws.Range("K110").Value = 1649
$ws.Range("L110").Value = 3701.8572
$ws.Range("M110").Value = 396
$ws.Range("N110").Value = -7791.8572
# Row 117
$ws.Range("H117").Value = 28146
$ws.Range("J117").Value = 28146
$ws.Range("L117").Value = 28146
$ws.Range("N117").Value = -37324

$ws = $wb.Worksheets.Item("BSM")
# Row 20
$ws.Range("H20").Value = 3046
$ws.Range("I20").Value = 2579.5
$ws.Range("J20").Value = 3512.5
$ws.Range("K20").Value = 2579.5
$ws.Range("L20").Value = 3512.5
$ws.Range("M20").Value = -2332.5
$ws.Range("N20").Value = -4006.5
# Row 27
$ws.Range("H27").Value = 0
$ws.Range("J27").Value = 0
$ws.Range("L27").Value = 0
$ws.Range("N27").ClearContents()
# Row 134
$ws.Range("H134").Value = 1388.849
$ws.Range("I134").Value = 1157.6459
$ws.Range("J134").Value = 3608.4
$ws.Range("K134").Value = 3472.9377
$ws.Range("L134").Value = 10825.2
$ws.Range("M134").Value = -937.9377
$ws.Range("N134").Value = -15895.2

$ws = $wb.Worksheets.Item("CRP")
# Row 33
$ws.Range("H33").Value = 22233
$ws.Range("I33").Value = 2020.6666
$ws.Range("J33").Value = 32339.166
$ws.Range("K33").Value = 2020.6666
$ws.Range("L33").Value = 32339.166
$ws.Range("M33").Value = -1641.6666
$ws.Range("N33").Value = -33097.166
# Row 36
$ws.Range("H36").Value = 18000
$ws.Range("J36").Value = 18000
$ws.Range("L36").Value = 18000
$ws.Range("N36").Value = -18776
# Row 40
$ws.Range("H40").Value = 18000
$ws.Range("J40").Value = 18000
$ws.Range("L40").Value = 18000
$ws.Range("N40").Value = -18320
# Row 42
$ws.Range("H42").Value = 17000
$ws.Range("J42").Value = 17000
$ws.Range("L42").Value = 17000
$ws.Range("N42").Value = -18186
# Row 44
$ws.Range("H44").Value = 34750
$ws.Range("J44").Value = 34750
$ws.Range("L44").Value = 34750
$ws.Range("N44").Value = -35634
# Row 107
$ws.Range("H107").Value = 689
$ws.Range("I107").Value = 325
$ws.Range("J107").Value = 1019.9091
$ws.Range("K107").Value = 325
$ws.Range("L107").Value = 1019.9091
$ws.Range("M107").Value = 1595
$ws.Range("N107").Value = -4859.9091
# Row 132
$ws.Range("H132").Value = 2117.8865
$ws.Range("I132").Value = 2066.919
$ws.Range("J132").Value = 2387.2856
$ws.Range("K132").Value = 6200.757
$ws.Range("L132").Value = 7161.8568
$ws.Range("M132").Value = -3670.757
$ws.Range("N132").Value = -12221.8568

$ws = $wb.Worksheets.Item("CUL")
# Row 114
$ws.Range("H114").Value = 2674.889
$ws.Range("I114").Value = 1611.7273
$ws.Range("J114").Value = 4345.5713
$ws.Range("K114").Value = 4835.1819
$ws.Range("L114").Value = 13036.7139
$ws.Range("M114").Value = -1581.1819
$ws.Range("N114").Value = -19544.7139
# Row 117
$ws.Range("H117").Value = 4108.909
$ws.Range("I117").Value = 3073.125
$ws.Range("J117").Value = 4700.7856
$ws.Range("K117").Value = 9219.375
$ws.Range("L117").Value = 14102.3568
$ws.Range("M117").Value = -5777.375
$ws.Range("N117").Value = -20986.3568
# Row 121
$ws.Range("H121").Value = 7143811
$ws.Range("I121").Value = 363.75
$ws.Range("J121").Value = 16668407
$ws.Range("K121").Value = 1091.25
$ws.Range("L121").Value = 50005221
$ws.Range("M121").Value = 218.75
$ws.Range("N121").Value = -50007841
# Row 140
$ws.Range("H140").Value = 12297.484
$ws.Range("I140").Value = 13289.519
$ws.Range("K140").Value = 39868.557
$ws.Range("M140").Value = -34688.557

$ws = $wb.Worksheets.Item("GSM")
# Row 70
$ws.Range("H70").Value = 7846.154
$ws.Range("J70").Value = 4800
$ws.Range("L70").Value = 4800
$ws.Range("N70").Value = -5340
# Row 73
$ws.Range("H73").Value = 7846.154
$ws.Range("J73").Value = 4800
$ws.Range("L73").Value = 4800
$ws.Range("N73").Value = -6672

$ws = $wb.Worksheets.Item("LTW")
# Row 47
$ws.Range("H47").Value = 77565
$ws.Range("J47").Value = 77565
$ws.Range("L47").Value = 77565
$ws.Range("N47").Value = -78545
# Row 52
$ws.Range("H52").Value = 77565
$ws.Range("J52").Value = 77565
$ws.Range("L52").Value = 77565
$ws.Range("N52").Value = -78031

$ws = $wb.Worksheets.Item("WVR")
# Row 113
$ws.Range("H113").Value = 442.66666
$ws.Range("I113").Value = 442.66666
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 1327.99998
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = 842.00002
$ws.Range("N113").ClearContents()
# Row 114
$ws.Range("H114").Value = 0
$ws.Range("J114").Value = 0
$ws.Range("L114").Value = 0
$ws.Range("N114").ClearContents()

Write-Host "applied leve-profit updates"
